$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "wenden"
$ws.Range("B2").Value = "none"
$ws.Range("C2").Value = "none"
$ws.Range("A3").Value = "landen"
$ws.Range("B3").Value = "flower/flower020.jpg"
$ws.Range("C3").Value = "flower"
$ws.Range("A4").Value = "steuern"
$ws.Range("B4").Value = "flower/flower014.jpg"
$ws.Range("C4").Value = "flower"
$ws.Range("A5").Value = "rufen"
$ws.Range("B5").Value = "none"
$ws.Range("C5").Value = "none"
$ws.Range("A6").Value = "passen"
$ws.Range("B6").Value = "dog/dog002.jpg"
$ws.Range("C6").Value = "dog"
$ws.Range("A7").Value = "geben"
$ws.Range("B7").Value = "dog/dog016.jpg"
$ws.Range("C7").Value = "dog"
$ws.Range("A8").Value = "enden"
$ws.Range("B8").Value = "none"
$ws.Range("C8").Value = "none"
$ws.Range("A9").Value = "bitten"
$ws.Range("B9").Value = "dog/dog022.jpg"
$ws.Range("C9").Value = "dog"
$ws.Range("A10").Value = "zahlen"
$ws.Range("B10").Value = "flower/flower010.jpg"
$ws.Range("C10").Value = "flower"
$ws.Range("A11").Value = "opfern"
$ws.Range("B11").Value = "none"
$ws.Range("C11").Value = "none"
$ws.Range("A12").Value = "hoffen"
$ws.Range("B12").Value = "dog/dog008.jpg"
$ws.Range("C12").Value = "dog"
$ws.Range("A13").Value = "heilen"
$ws.Range("B13").Value = "flower/flower000.jpg"
$ws.Range("C13").Value = "flower"
$ws.Range("A14").Value = "schalten"
$ws.Range("B14").Value = "none"
$ws.Range("C14").Value = "none"
$ws.Range("A15").Value = "schulden"
$ws.Range("B15").Value = "flower/flower030.jpg"
$ws.Range("C15").Value = "flower"
$ws.Range("A16").Value = "leuchten"
$ws.Range("B16").Value = "flower/flower023.jpg"
$ws.Range("C16").Value = "flower"
$ws.Range("A17").Value = "klagen"
$ws.Range("B17").Value = "none"
$ws.Range("C17").Value = "none"
$ws.Range("A18").Value = "faulen"
$ws.Range("B18").Value = "flower/flower018.jpg"
$ws.Range("C18").Value = "flower"
$ws.Range("A19").Value = "lächeln"
$ws.Range("B19").Value = "flower/flower019.jpg"
$ws.Range("C19").Value = "flower"
$ws.Range("A20").Value = "sparen"
$ws.Range("B20").Value = "none"
$ws.Range("C20").Value = "none"
$ws.Range("A21").Value = "runden"
$ws.Range("B21").Value = "dog/dog000.jpg"
$ws.Range("C21").Value = "dog"
$ws.Range("A22").Value = "buchen"
$ws.Range("B22").Value = "dog/dog004.jpg"
$ws.Range("C22").Value = "dog"
$ws.Range("A23").Value = "hören"
$ws.Range("B23").Value = "none"
$ws.Range("C23").Value = "none"
$ws.Range("A24").Value = "stärken"
$ws.Range("B24").Value = "dog/dog019.jpg"
$ws.Range("C24").Value = "dog"
$ws.Range("A25").Value = "zeugen"
$ws.Range("B25").Value = "dog/dog006.jpg"
$ws.Range("C25").Value = "dog"
$ws.Range("A26").Value = "drohen"
$ws.Range("B26").Value = "none"
$ws.Range("C26").Value = "none"
$ws.Range("A27").Value = "küssen"
$ws.Range("B27").Value = "flower/flower026.jpg"
$ws.Range("C27").Value = "flower"
$ws.Range("A28").Value = "platzen"
$ws.Range("B28").Value = "flower/flower027.jpg"
$ws.Range("C28").Value = "flower"
$ws.Range("A29").Value = "orten"
$ws.Range("B29").Value = "none"
$ws.Range("C29").Value = "none"
$ws.Range("A30").Value = "schütteln"
$ws.Range("B30").Value = "dog/dog020.jpg"
$ws.Range("C30").Value = "dog"
$ws.Range("A31").Value = "segnen"
$ws.Range("B31").Value = "flower/flower008.jpg"
$ws.Range("C31").Value = "flower"
$ws.Range("A32").Value = "weigern"
$ws.Range("B32").Value = "none"
$ws.Range("C32").Value = "none"
$ws.Range("A33").Value = "dürfen"
$ws.Range("B33").Value = "flower/flower002.jpg"
$ws.Range("C33").Value = "flower"
$ws.Range("A34").Value = "dringen"
$ws.Range("B34").Value = "dog/dog031.jpg"
$ws.Range("C34").Value = "dog"
$ws.Range("A35").Value = "ächzen"
$ws.Range("B35").Value = "none"
$ws.Range("C35").Value = "none"
$ws.Range("A36").Value = "parken"
$ws.Range("B36").Value = "dog/dog028.jpg"
$ws.Range("C36").Value = "dog"
$ws.Range("A37").Value = "starren"
$ws.Range("B37").Value = "flower/flower017.jpg"
$ws.Range("C37").Value = "flower"
$ws.Range("A38").Value = "kosten"
$ws.Range("B38").Value = "none"
$ws.Range("C38").Value = "none"
$ws.Range("A39").Value = "leisten"
$ws.Range("B39").Value = "dog/dog013.jpg"
$ws.Range("C39").Value = "dog"
$ws.Range("A40").Value = "streichen"
$ws.Range("B40").Value = "dog/dog026.jpg"
$ws.Range("C40").Value = "dog"
$ws.Range("A41").Value = "dauern"
$ws.Range("B41").Value = "none"
$ws.Range("C41").Value = "none"
$ws.Range("A42").Value = "werden"
$ws.Range("B42").Value = "flower/flower006.jpg"
$ws.Range("C42").Value = "flower"
$ws.Range("A43").Value = "heben"
$ws.Range("B43").Value = "flower/flower013.jpg"
$ws.Range("C43").Value = "flower"
$ws.Range("A44").Value = "stören"
$ws.Range("B44").Value = "none"
$ws.Range("C44").Value = "none"
$ws.Range("A45").Value = "schützen"
$ws.Range("B45").Value = "flower/flower029.jpg"
$ws.Range("C45").Value = "flower"
$ws.Range("A46").Value = "fühlen"
$ws.Range("B46").Value = "dog/dog018.jpg"
$ws.Range("C46").Value = "dog"
$ws.Range("A47").Value = "bremsen"
$ws.Range("B47").Value = "none"
$ws.Range("C47").Value = "none"
$ws.Range("A48").Value = "schultern"
$ws.Range("B48").Value = "dog/dog014.jpg"
$ws.Range("C48").Value = "dog"
$ws.Range("A49").Value = "mauern"
$ws.Range("B49").Value = "dog/dog003.jpg"
$ws.Range("C49").Value = "dog"
